$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the monthly sheets: strip the accented "ê" and fix "Março" typo.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Mês - Janeiro").Name   = "Mes - Janeiro"
$wb.Worksheets.Item("Mês - Fevereiro").Name = "Mes - Fevereiro"
$wb.Worksheets.Item("Mês - Março").Name     = "Mes - Marco"
$wb.Worksheets.Item("Mês - Abril").Name     = "Mes - Abril"
$wb.Worksheets.Item("Mês - Maio").Name      = "Mes - Maio"
$wb.Worksheets.Item("Mês - Julho").Name     = "Mes - Julho"

# ---------------------------------------------------------------------------
# 2) "Geral" sheet: fill in the partial totals for round 38 (column B) and
#    append a new "Lider_Rodada" row at the bottom.
# ---------------------------------------------------------------------------
$geral = $wb.Worksheets.Item("Geral")

$geralBValues = @(59.56, 45.3, 30.06, 69.56, 54.06, 71.70999999999999, 53.06, 69.26000000000001, 71.36, 47.2, 58.4, 45.46, 48.46, 30.32, 52.36, 50.6, 44.06, 49.16, 21.7, 56.65, 43.1, 41.6, 39.66, 39.66, 54.36, 58.26, 58.51, 69.76000000000001, 47.66, 45.86, 23.26, 73.66, 60.36, 26.51, 48.06, 60.66, 51.05, 59.76, 34.36, 38.5, 55.96, 34.76, 45.3, 57.76, 50.1, 64.95999999999999, 50.76, 61.8, 59.86, 46.56, 63.66, 35.8)

for ($i = 0; $i -lt $geralBValues.Length; $i++) {
    $row = $i + 2
    $geral.Cells.Item($row, 2).Value = $geralBValues[$i]
}

# New row 54 - "Lider_Rodada" team, all zeros (A:AM = columns 1:39).
$geral.Cells.Item(54, 1).Value = "Lider_Rodada"
$geral.Range("A53").Copy()
$geral.Cells.Item(54, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
for ($c = 2; $c -le 39; $c++) {
    $geral.Cells.Item(54, $c).Value = 0
}

# ---------------------------------------------------------------------------
# 3) "Turno 2" sheet: append the same new "Lider_Rodada" row (columns A:T).
# ---------------------------------------------------------------------------
$turno2 = $wb.Worksheets.Item("Turno 2")

$turno2.Cells.Item(54, 1).Value = "Lider_Rodada"
$turno2.Range("A53").Copy()
$turno2.Cells.Item(54, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
for ($c = 2; $c -le 20; $c++) {
    $turno2.Cells.Item(54, $c).Value = 0
}

# ---------------------------------------------------------------------------
# 4) "Classif Turno 2" sheet: re-rank rows 2:53 by the (now non-zero) round
#    38 totals, descending - mirroring the Geral sheet's own row order.
# ---------------------------------------------------------------------------
$classif = $wb.Worksheets.Item("Classif Turno 2")

# For each destination row (2..53) in "Classif Turno 2", the matching source
# row in "Geral" that holds that rank's team/value.
$classifSrcRow = @(33, 7, 10, 29, 5, 9, 47, 52, 49, 37, 34, 50, 39, 2, 28, 12, 27, 45, 21, 42, 26, 6, 8, 16, 38, 48, 17, 46, 19, 14, 36, 30, 11, 51, 31, 13, 44, 3, 18, 22, 23, 24, 25, 41, 53, 43, 40, 15, 4, 35, 32, 20)

for ($i = 0; $i -lt $classifSrcRow.Length; $i++) {
    $destRow = $i + 2
    $srcRow = $classifSrcRow[$i]
    $classif.Cells.Item($destRow, 1).Value = $geral.Cells.Item($srcRow, 1).Value2
    $classif.Cells.Item($destRow, 2).Value = $geral.Cells.Item($srcRow, 2).Value2
}
